$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the D2:D22 "WolframAlpha Value" data and D1 header, keeping formatting
$ws.Range("D1").ClearContents()
$ws.Range("D2:D22").ClearContents()

# Update selection to match the new state
$ws.Range("D1:D22").Select()
